$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column B ("CUPOS") - pushes the old "DIRECCIÓN" column to C
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).Insert()

# ---------------------------------------------------------------------------
# 2. Header text + CUPOS values
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "CUPOS"

$values = @(4,4,4,4,4,3,3,3,3,3,3,2,2,2,2,2,2,1,1,1,1,1,1,1,1,1,1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# ---------------------------------------------------------------------------
# 3. Font overhaul: Georgia -> Arial across the whole table, headers bumped
#    to size 11
# ---------------------------------------------------------------------------
# Body text columns (A institution names, C addresses) keep size 10
$ws.Range("A2:A28").Font.Name = "Arial"
$ws.Range("C2:C28").Font.Name = "Arial"

# Header row (A1, C1): bold, Arial, size 11, explicit black color
$ws.Range("A1").Font.Name = "Arial"
$ws.Range("A1").Font.Size = 11
$ws.Range("A1").Font.Bold = $true

$ws.Range("C1").Font.Name = "Arial"
$ws.Range("C1").Font.Size = 11
$ws.Range("C1").Font.Bold = $true

# New "CUPOS" header (B1): bold, Arial, size 11, automatic (theme) color,
# no border / no fill, centered, no wrap
$ws.Range("B1").Font.Name = "Arial"
$ws.Range("B1").Font.Size = 11
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Borders.LineStyle = -4142
$ws.Range("B1").Interior.Pattern = -4142
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4108
$ws.Range("B1").WrapText = $false

# New "CUPOS" values (B2:B28): plain, Arial, size 11, automatic color,
# bordered, centered, no wrap, no fill
$ws.Range("B2:B28").Font.Name = "Arial"
$ws.Range("B2:B28").Font.Size = 11
$ws.Range("B2:B28").Interior.Pattern = -4142
$ws.Range("B2:B28").HorizontalAlignment = -4108
$ws.Range("B2:B28").VerticalAlignment = -4108
$ws.Range("B2:B28").WrapText = $false

# ---------------------------------------------------------------------------
# 4. New (blank) trailing row 29
# ---------------------------------------------------------------------------
$ws.Range("A29:C29").Value = ""

# ---------------------------------------------------------------------------
# 5. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 45.7109375
$ws.Columns.Item(2).ColumnWidth = 30
$ws.Columns.Item(3).ColumnWidth = 50.140625

# ---------------------------------------------------------------------------
# 6. Row heights adjustments (match the authored sizing)
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 26.25

# ---------------------------------------------------------------------------
# 7. Active cell / selection ends on C1 after the edit
# ---------------------------------------------------------------------------
$ws.Range("C1").Select()

Write-Output "done"
